$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2-5 (values changed, text unchanged)
$ws.Range("B2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 63
$ws.Range("F2").Value = 4899
$ws.Range("G2").Value = 1.285976729944887

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 26
$ws.Range("G3").Value = 1.024428684003152

$ws.Range("B4").Value = 0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 4
$ws.Range("G4").Value = 0.04908577739599951

$ws.Range("B5").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4
$ws.Range("G5").Value = 0.154619250096637

# Insert a new row 6 for "wildfly" before the old row 6 ("storm"), shifting storm to row 7
$ws.Rows.Item(6).Insert()

# New row 6: wildfly
$ws.Range("A6").Value = "wildfly"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 14079
$ws.Range("G6").Value = 0.01420555437175936

# Row 7: storm (updated values)
$ws.Range("A7").Value = "storm"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 3398
$ws.Range("G7").Value = 0.05885815185403178
